# Remove open burning and natural emissions
#
# The "map" sheet's rows for the open-burning / natural-emissions inventory
# sectors (4F Agricultural-residue-burning, 5A Forest-fires, 5C / 5D
# Other-natural) should no longer be scaled against a CEDS sector. We clear
# the "scaling_sector" (column B) value for those rows and instead flag them
# in a new "Notes" column (D) explaining they're excluded from the CMIP6
# data product.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# New "Notes" column header
$ws.Range("D1").Value = "Notes"

$note = "Don't include in CMIP6 data product"

# 4F_... (row 35), 5A_... (row 36), 5C_... (row 37), 5D_... (row 38):
# drop the scaling_sector mapping and annotate with the note.
$ws.Range("B35").ClearContents()
$ws.Range("D35").Value = $note

$ws.Range("B36").ClearContents()
$ws.Range("D36").Value = $note

$ws.Range("B37").ClearContents()
$ws.Range("D37").Value = $note

$ws.Range("B38").ClearContents()
$ws.Range("D38").Value = $note

# Leave the selection where the author ended up after the edit.
$ws.Activate() | Out-Null
$ws.Range("B39").Select() | Out-Null
